# "skipp comment and new example"
# Replace the bank-payment-details example (with the stray "Payee's bank"
# comment line) with a small generic numbered example that demonstrates
# merged cells (1..5 laid out over six merge regions on a 5x5 grid).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- start clean -----------------------------------------------------
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# --- borders for the whole 5x5 table ---------------------------------
$all = $ws.Range("A1:E5")
$all.Borders.LineStyle = 1

# --- force the numeric-looking labels to be stored as TEXT ------------
# (so "1".."5" land in sharedStrings instead of becoming numbers)
$all.NumberFormat = "@"

# --- cell values -------------------------------------------------------
$ws.Range("A1").Value = "1"
$ws.Range("D1").Value = "4"
$ws.Range("E1").Value = "5`n"

$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "2`n"
$ws.Range("C2").Value = "3"
$ws.Range("D2").Value = "4"

$ws.Range("A3").Value = "1`n"
$ws.Range("C3").Value = "3"
$ws.Range("D3").Value = "4"

$ws.Range("B4").Value = "2"
$ws.Range("E4").Value = "5"

$ws.Range("A5").Value = "1"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "3"

# --- wrap + top/left alignment on the three merge regions that hold
#     the multi-line labels -------------------------------------------
foreach ($addr in @("E1:E3", "B2:B3", "A3:A4")) {
    $r = $ws.Range($addr)
    $r.WrapText = $true
    $r.HorizontalAlignment = -4131
    $r.VerticalAlignment = -4160
}

# --- merges --------------------------------------------------------
$ws.Range("A1:C1").Merge()
$ws.Range("E1:E3").Merge()
$ws.Range("B2:B3").Merge()
$ws.Range("A3:A4").Merge()
$ws.Range("B4:D4").Merge()
$ws.Range("C5:E5").Merge()

Write-Output "done"
